$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# row -> @(nombre_aides, montant_total)
$changes = @{
    45  = @(23378, 99234004)
    47  = @(3598, 31499398)
    54  = @(60294, 353615499)
    61  = @(51959, 133917851)
    92  = @(409173, 1595662657)
    94  = @(94212, 918364469)
    95  = @(50780, 933227656)
    97  = @(2161, 214282109)
    104 = @(135248, 272240066)
    141 = @(80475, 280728664)
    142 = @(168976, 681799746)
    154 = @(201571, 786779728)
    182 = @(71, 11214004)
}

foreach ($row in $changes.Keys) {
    $vals = $changes[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
}
